$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "840"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1887687.79"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1032"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3681879.47"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "670"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2190407.78"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "377"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1422220.18"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "290"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "749727.45"

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "573"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2392474.89"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "396"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1369778.57"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "400"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1073299.43"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "633"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2543733.99"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "438"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1548161.40"

$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "389"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "985306.70"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "526"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "1735225.47"

$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "38"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "104500.00"

$ws.Range("C92").NumberFormat = "@"
$ws.Range("C92").Value = "669"
$ws.Range("D92").NumberFormat = "@"
$ws.Range("D92").Value = "1632584.94"

